# Update Name of Algo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.068
$ws.Range("D3").Value = -7.473000000000001
$ws.Range("E8").Value = 16.53
$ws.Range("E11").Value = 16.833
$ws.Range("A12").Value = -21.544
$ws.Range("B14").Value = 5.981
$ws.Range("E14").Value = 16.887
$ws.Range("E15").Value = 16.223
$ws.Range("B26").Value = 6.002
$ws.Range("D30").Value = -7.257
$ws.Range("B31").Value = 6.29
$ws.Range("A32").Value = -21.345
$ws.Range("B35").Value = 8.095000000000001
$ws.Range("A36").Value = -21.147
$ws.Range("E36").Value = 16.263
$ws.Range("B37").Value = 8.260000000000002
$ws.Range("A38").Value = -20.093
$ws.Range("D44").Value = -7.874
$ws.Range("B45").Value = 5.933000000000001
$ws.Range("A46").Value = -21.553
$ws.Range("A54").Value = -22.209
$ws.Range("A55").Value = -22.21
$ws.Range("B57").Value = 5.331999999999999
$ws.Range("D58").Value = -8.303000000000001
$ws.Range("E64").Value = 17.185
$ws.Range("A67").Value = -21.603
$ws.Range("A69").Value = -21.636
$ws.Range("A72").Value = -21.55
$ws.Range("D84").Value = -8.216999999999999
$ws.Range("D89").Value = -7.105
$ws.Range("E89").Value = 17.149
$ws.Range("A91").Value = -21.608
$ws.Range("D91").Value = -6.910000000000001
$ws.Range("D92").Value = -6.736
$ws.Range("A99").Value = -20.828
$ws.Range("B100").Value = 5.558
$ws.Range("B102").Value = 7.468000000000001
$ws.Range("D102").Value = -7.873
